# Release-Notes.xlsx update: the folder inventory scrape no longer contains the
# "Build-Custom-Knowledge-RAG-App-With-Azure-AI-Foundry" folder (it dropped out
# of the latest scrape), so its row is removed from the Folder Inventory sheet
# and every row below it shifts up by one. The Metadata and Summary sheets are
# refreshed to match the new inventory: the "Generated On" timestamp and
# "Workflow Run" counter move forward, and the folder-count totals drop by one.

$wb = $excel.ActiveWorkbook

# --- Folder Inventory: remove the row for the folder that disappeared ---
$inventory = $wb.Worksheets.Item("Folder Inventory")
$inventory.Rows.Item(3).Delete()

# --- Metadata sheet: refresh generated-on timestamp, folder count, run id ---
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B3").Value = "2025-06-13 07:47:12 UTC"
$metadata.Range("B4").Value = 74

# "Workflow Run" is stored as text in the workbook (not a number), so force
# a text number format before assigning the numeric-looking string — otherwise
# Excel would silently coerce it to a numeric cell.
$metadata.Range("B5").NumberFormat = "@"
$metadata.Range("B5").Value = "20"
$metadata.Range("B5").Style = "Normal"

# --- Summary sheet: refresh total folder counters ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 74
$summary.Range("B3").Value = 74
